$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update 想去人数 (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1429
$wsExpo.Range("F3").Value = 3006
$wsExpo.Range("F4").Value = 33
$wsExpo.Range("F5").Value = 194
$wsExpo.Range("F6").Value = 280

# Sheet "全部类型" (All types) - update 想去人数 (want-to-go count) column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1429
$wsAll.Range("F3").Value = 3006
$wsAll.Range("F4").Value = 33
$wsAll.Range("F5").Value = 194
$wsAll.Range("F7").Value = 280
